# The deck's theme was re-pointed: the Slide Master (ppt/theme/theme1.xml,
# previously the "Integral" theme) is switched to use the built-in
# "Office Theme" palette. (The font scheme and format scheme are identical
# between the two themes - only the 12 color-scheme slots differ.)
#
# PowerPoint's object model exposes those 12 DrawingML theme colors
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) through the classic
# ColorScheme/Colors(1..12).RGB properties, in that fixed slot order.

function Convert-RGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.ColorScheme

# Target palette: the built-in "Office Theme" colors.
$officeTheme = @(
    @(0x00, 0x00, 0x00),  # 1  dk1
    @(0xFF, 0xFF, 0xFF),  # 2  lt1
    @(0x44, 0x54, 0x6A),  # 3  dk2
    @(0xE7, 0xE6, 0xE6),  # 4  lt2
    @(0x5B, 0x9B, 0xD5),  # 5  accent1
    @(0xED, 0x7D, 0x31),  # 6  accent2
    @(0xA5, 0xA5, 0xA5),  # 7  accent3
    @(0xFF, 0xC0, 0x00),  # 8  accent4
    @(0x44, 0x72, 0xC4),  # 9  accent5
    @(0x70, 0xAD, 0x47),  # 10 accent6
    @(0x05, 0x63, 0xC1),  # 11 hlink
    @(0x95, 0x4F, 0x72)   # 12 folHlink
)

for ($i = 0; $i -lt $officeTheme.Length; $i++) {
    $rgbTriplet = $officeTheme[$i]
    $colorScheme.Colors($i + 1).RGB = Convert-RGB $rgbTriplet[0] $rgbTriplet[1] $rgbTriplet[2]
}
